$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Invoice Number"
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("A2:A11").Select()
